$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

$ws.Range("D3").Value = "DateOfBirth`n`tتاريخ الميلاد`nDD/MM/YYYY"
$ws.Range("E3").Value = "DateOfBirthHijri`n`tتاريخ الميلاد هجري`nYYYYMMDD"

$colH = $tbl.ListColumns.Add()
$colI = $tbl.ListColumns.Add()

$ws.Range("G3").Copy()
$ws.Range("H3").PasteSpecial(-4122)
$ws.Range("F3").Copy()
$ws.Range("I3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H3").Value = "SpecialNeed`t`nالإحتياجات الخاصة`nYes/No"
$ws.Range("I3").Value = "SpecialNeedList`n`tقائمة الإحتياجات الخاصة"

# body rows 4-30 (not last row 31) : copy from G (full box)
$ws.Range("G4:G30").Copy()
$ws.Range("H4:H30").PasteSpecial(-4122)
$ws.Range("I4:I30").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# last row 31: copy from G31 (box minus bottom)
$ws.Range("G31").Copy()
$ws.Range("H31").PasteSpecial(-4122)
$ws.Range("I31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# extra rows below the table (32-33 body style, 34 last-row style), columns H:I only
$ws.Range("H4:I4").Copy()
$ws.Range("H32:I33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H31:I31").Copy()
$ws.Range("H34:I34").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# column widths for the two new columns
$ws.Range("H1:I1").ColumnWidth = 20.28

# leave selection on the new header cell, matching the saved view
$ws.Range("H3").Select()
